# fix(publipostage): Correct status name
#
# Replace the "bleu" status label with "noir", and reword the
# statut_name strings describing results/publications, across the
# relevant cells of the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (statut_label): "bleu" -> "noir"
foreach ($cell in @("B6", "B7", "B8")) {
    $ws.Range($cell).Value = "noir"
}

# Column C (statut_name): reword wording for each status
$ws.Range("C2").Value = "résultat postés ou publiés dans les 12 mois"
$ws.Range("C5").Value = "résultat postés ou publiés dans les 12 mois"

$ws.Range("C3").Value = "résultat postés ou publiés dans les 36 mois"

$ws.Range("C4").Value = "résultat postés ou publiés"

foreach ($cell in @("C6", "C7", "C8")) {
    $ws.Range($cell).Value = "pas de résultat postés ni publiés"
}
